# Burndown Chart update:
#  - Sheet1!C9 (Actual, day 8) drops from 30 to 19 (stretch-goal story points
#    completed), reflecting currently completed User Stories.
#  - Active selection moves to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Actual (column C) story points completed on day 8 (row 9) revised from 30 to 19
$ws.Range("C9").Value = 19

# Move the active cell selection to D9
$ws.Range("D9").Select() | Out-Null
